# Update Jeremy en Ivar
# Adds new timesheet rows to "P1 - Jeremy" and "P4 - Ivar" sheets, and
# nudges the saved selection/view state to match where the author ended up.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# P1 - Jeremy: add rows 20-24 with new activities
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("P1 - Jeremy")

$ws1.Range("A20").Value = "les KBS"
$ws1.Range("B20").Value = [DateTime]"2020-10-07"
$ws1.Range("C20").Value = 60

$ws1.Range("A21").Value = "Gastcolleges bekeken en aantekingen gemaakt"
$ws1.Range("B21").Value = [DateTime]"2020-10-08"
$ws1.Range("C21").Value = 60

$ws1.Range("A22").Value = "Opzet conversie verhogende maatregelen"
$ws1.Range("B22").Value = [DateTime]"2020-10-08"
$ws1.Range("C22").Value = 45

$ws1.Range("A23").Value = "Verder werken aan conversiemaatregelen"
$ws1.Range("B23").Value = [DateTime]"2020-10-12"
$ws1.Range("C23").Value = 105

$ws1.Range("A24").Value = "Afmaken conversiemaatregelen"
$ws1.Range("B24").Value = [DateTime]"2020-10-13"
$ws1.Range("C24").Value = 50

[void]$ws1.Range("D24").Select()

# ---------------------------------------------------------------------------
# P4 - Ivar: add rows 30-31 with new activities
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("P4 - Ivar")

$ws4.Range("A30").Value = "Voortzet conversie verhogende maatregelen"
$ws4.Range("B30").Value = [DateTime]"2020-10-13"
$ws4.Range("C30").Value = 50

$ws4.Range("A31").Value = "Weekverslag ingevuld"
$ws4.Range("B31").Value = [DateTime]"2020-10-08"
$ws4.Range("C31").Value = 5

# Row 18 re-wraps slightly shorter after the sheet layout refresh.
$ws4.Rows.Item(18).RowHeight = 43.2

# ---------------------------------------------------------------------------
# P6 - Jasper: just a leftover selection change from the author's session
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("P6 - Jasper")
[void]$ws6.Range("C31").Select()

# Re-select P4 - Ivar so it stays the active sheet/cell as in the saved file.
[void]$ws4.Activate()
[void]$ws4.Range("D31").Select()
